$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2951.1428
$ws.Range("J17").Value = 3137.3845
$ws.Range("L17").Value = 9412.1535
$ws.Range("N17").Value = -9748.1535

$ws.Range("H40").Value = 1322.3572
$ws.Range("I40").Value = 1580
$ws.Range("J40").Value = 1179.2222
$ws.Range("K40").Value = 1580
$ws.Range("L40").Value = 1179.2222
$ws.Range("M40").Value = -1405
$ws.Range("N40").Value = -1529.2222

$ws.Range("H41").Value = 416.8889
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 813
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 813
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -1693

$ws.Range("H82").Value = 676.8
$ws.Range("I82").Value = 676.8
$ws.Range("K82").Value = 2030.4
$ws.Range("M82").Value = -1624.4

$ws.Range("H85").Value = 676.8
$ws.Range("I85").Value = 676.8
$ws.Range("K85").Value = 2030.4
$ws.Range("M85").Value = -626.3999999999999

$ws.Range("H138").Value = 5181.0137
$ws.Range("I138").Value = 1350.7142
$ws.Range("J138").Value = 7618.477
$ws.Range("K138").Value = 4052.1426
$ws.Range("L138").Value = 22855.431
$ws.Range("M138").Value = 1087.8574
$ws.Range("N138").Value = -33135.431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1549.5385
$ws.Range("I2").Value = 1569.0435
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1569.0435
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1456.0435
$ws.Range("N2").Value = -1626

$ws.Range("H45").Value = 1537.8182
$ws.Range("I45").Value = 1491.1333
$ws.Range("J45").Value = 2004.6666
$ws.Range("K45").Value = 1491.1333
$ws.Range("L45").Value = 2004.6666
$ws.Range("M45").Value = -1114.1333
$ws.Range("N45").Value = -2758.6666

$ws.Range("H102").Value = 3460.7273
$ws.Range("I102").Value = 2811.5
$ws.Range("K102").Value = 2811.5
$ws.Range("M102").Value = -1189.5

$ws.Range("H116").Value = 1549.5385
$ws.Range("I116").Value = 1569.0435
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1569.0435
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 724.9565
$ws.Range("N116").Value = -5988

$ws.Range("H132").Value = 2504.6
$ws.Range("I132").Value = 1693.4706
$ws.Range("J132").Value = 4228.25
$ws.Range("K132").Value = 5080.4118
$ws.Range("L132").Value = 12684.75
$ws.Range("M132").Value = -2550.4118
$ws.Range("N132").Value = -17744.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1549.5385
$ws.Range("I3").Value = 1569.0435
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1569.0435
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1455.0435
$ws.Range("N3").Value = -1628

$ws.Range("H20").Value = 1031.7778
$ws.Range("I20").Value = 1302
$ws.Range("J20").Value = 815.6
$ws.Range("K20").Value = 1302
$ws.Range("L20").Value = 815.6
$ws.Range("M20").Value = -1055
$ws.Range("N20").Value = -1309.6

$ws.Range("H107").Value = 1505.9
$ws.Range("I107").Value = 1008.25
$ws.Range("J107").Value = 1837.6666
$ws.Range("K107").Value = 1008.25
$ws.Range("L107").Value = 1837.6666
$ws.Range("M107").Value = 911.75
$ws.Range("N107").Value = -5677.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7737.8
$ws.Range("I16").Value = 9213.916999999999
$ws.Range("J16").Value = 1833.3334
$ws.Range("K16").Value = 9213.916999999999
$ws.Range("L16").Value = 1833.3334
$ws.Range("M16").Value = -8926.916999999999
$ws.Range("N16").Value = -2407.3334

$ws.Range("H22").Value = 173.33333
$ws.Range("I22").Value = 165
$ws.Range("J22").Value = 190
$ws.Range("K22").Value = 165
$ws.Range("L22").Value = 190
$ws.Range("M22").Value = 185
$ws.Range("N22").Value = -890

$ws.Range("H28").Value = 34000
$ws.Range("J28").Value = 34000
$ws.Range("L28").Value = 34000
$ws.Range("N28").Value = -34490

$ws.Range("H113").Value = 7737.8
$ws.Range("I113").Value = 9213.916999999999
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 9213.916999999999
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = -7043.916999999999
$ws.Range("N113").Value = -6173.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 354.64706
$ws.Range("I15").Value = 216.35715
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 649.0714499999999
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = -509.0714499999999
$ws.Range("N15").Value = -3280

$ws.Range("H131").Value = 28573.229
$ws.Range("J131").Value = 81033.086
$ws.Range("L131").Value = 243099.258
$ws.Range("N131").Value = -253179.258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2626.0454
$ws.Range("I132").Value = 2365.9375
$ws.Range("J132").Value = 3319.6667
$ws.Range("K132").Value = 7097.8125
$ws.Range("L132").Value = 9959.000100000001
$ws.Range("M132").Value = -4567.8125
$ws.Range("N132").Value = -15019.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2127.7778
$ws.Range("I82").Value = 1725
$ws.Range("J82").Value = 2933.3333
$ws.Range("K82").Value = 1725
$ws.Range("L82").Value = 2933.3333
$ws.Range("M82").Value = -1364
$ws.Range("N82").Value = -3655.3333

$ws.Range("H85").Value = 2127.7778
$ws.Range("I85").Value = 1725
$ws.Range("J85").Value = 2933.3333
$ws.Range("K85").Value = 1725
$ws.Range("L85").Value = 2933.3333
$ws.Range("M85").Value = -477
$ws.Range("N85").Value = -5429.3333

$ws.Range("H122").Value = 6103.224
$ws.Range("I122").Value = 5387.05
$ws.Range("J122").Value = 7694.722
$ws.Range("K122").Value = 16161.15
$ws.Range("L122").Value = 23084.166
$ws.Range("M122").Value = -13711.15
$ws.Range("N122").Value = -27984.166

$ws.Range("H132").Value = 4580.979
$ws.Range("I132").Value = 6136.077
$ws.Range("J132").Value = 2743.1365
$ws.Range("K132").Value = 18408.231
$ws.Range("L132").Value = 8229.4095
$ws.Range("M132").Value = -15878.231
$ws.Range("N132").Value = -13289.4095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2059.353
$ws.Range("I122").Value = 1538.7693
$ws.Range("J122").Value = 3751.25
$ws.Range("K122").Value = 4616.3079
$ws.Range("L122").Value = 11253.75
$ws.Range("M122").Value = -2166.3079
$ws.Range("N122").Value = -16153.75
